# Delete row 324 ("止まりなさい" entry) from the posts sheet.
# This shifts all subsequent rows up by one, matching the target diff
# (dimension shrinks from A1:C483 to A1:C482).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(324).Delete()
